$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CEDEN_Benthic_Data_Dictionary")

# Insert two new rows before row 65. This pushes the existing "DataQuality" /
# "DataQualityIndicator" rows (old 65/66) down to 67/68, and all the trailing
# blank rows shift down by two as well.
$ws.Rows.Item(65).Insert()
$ws.Rows.Item(65).Insert()

# The freshly inserted rows pick up a "no border" flavor of the style used by
# the data rows above/below; copy the real row formatting down onto them so
# they match the rest of the table (border s=4 / s=5 / s=2).
$ws.Range("A64:F64").Copy()
$ws.Range("A65:F66").PasteSpecial(-4122)  # xlPasteFormats

# New row 65: ProgramCode
$ws.Range("A65").Value = "ProgramCode"
$ws.Range("B65").Value = "text"
$ws.Range("C65").FormulaArray = '=IFERROR(INDEX(Data_Dictionary_FromPDF!B:B,E65),"-")'
$ws.Range("D65").FormulaArray = '=IFERROR(INDEX(Data_Dictionary_FromPDF!D:D,E65),"-")'
$ws.Range("E65").Formula = "=MATCH(A65,Data_Dictionary_FromPDF!A:A,0)"
$ws.Range("F65").Value = "text"

# New row 66: ParentProjectCode
$ws.Range("A66").Value = "ParentProjectCode"
$ws.Range("B66").Value = "text"
$ws.Range("C66").FormulaArray = '=IFERROR(INDEX(Data_Dictionary_FromPDF!B:B,E66),"-")'
$ws.Range("D66").FormulaArray = '=IFERROR(INDEX(Data_Dictionary_FromPDF!D:D,E66),"-")'
$ws.Range("E66").Formula = "=MATCH(A66,Data_Dictionary_FromPDF!A:A,0)"
$ws.Range("F66").Value = "text"

# New column G notes for a few existing rows (EffortQACode, BenthicLabEffortComments,
# PercentSampleCounted).
$ws.Range("G52").Value = "Unique code applied to the result which describes any special conditions, situations or outliers occurring during or prior to lab sorting. Default value equals NR if unknown."
$ws.Range("G53").Value = "Comments related to lab sorting or sample processing."
$ws.Range("G54").Value = "Refers to the percent of the sample that was counted."
